$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Sales Order Identifier" (S2) and "Previous Doc" (AX2) both hold the
# same reconciled order number as text (it must not turn into a number,
# since leading/insignificant digits matter for this identifier). Force
# a text number format before writing the value so Excel keeps it as a
# string instead of auto-coercing the digit string into a number.
$ws.Range("S2").NumberFormat = "@"
$ws.Range("S2").Value = "3044924556"

$ws.Range("AX2").NumberFormat = "@"
$ws.Range("AX2").Value = "3044924556"
